# sub_plots in dashboard and minor improvments
#
# 1. Add a new "TORSO" worksheet at the end of the workbook with marker
#    data for the OSLO / CALGARY / JOHNATAN sites (mirrors the layout of
#    the other marker sheets: HIP, KNEE, ANKLE, FOOT, MTP, DISTAL_Marker).
# 2. Fix up the stale selection left on the MTP sheet (was pointing past
#    the used range at E8).
# 3. Leave the new TORSO sheet as the active / selected tab, like the
#    author did after finishing data entry on it.

$wb = $excel.ActiveWorkbook

# --- fix the MTP sheet's stale selection (E8 -> A1:E3) ------------------
$mtp = $wb.Worksheets.Item("MTP")
$mtp.Activate() | Out-Null
$mtp.Range("A1:E3").Select() | Out-Null

# --- add the new TORSO sheet after the last existing sheet --------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$torso = $wb.Worksheets.Add($null, $lastSheet)
$torso.Name = "TORSO"

# Header / site rows, same pattern as the other marker sheets.
$torso.Range("A1").Value = "OSLO"
$torso.Range("B1").Value = "C_7"
$torso.Range("C1").Value = "B_10"
$torso.Range("D1").Value = "sternum"
$torso.Range("E1").Value = "clav"

$torso.Range("A2").Value = "CALGARY"
$torso.Range("B2").Value = "C7"
$torso.Range("C2").Value = "T6"
$torso.Range("E2").Value = "STERNUM"
$torso.Range("D2").Value = "XYPHOID"

$torso.Range("A3").Value = "JOHNATAN"

# Leave the cursor where the author left it and make TORSO the active tab.
$torso.Range("D6").Select() | Out-Null
$torso.Activate() | Out-Null
